$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: new "tbd" (to-be-deleted) user, full sync-log line ---
$ws.Range("A21").Value = "tbd"
$ws.Range("B21").Value = "tbd"

$ws.Range("C21").Value = "tbd@eduvaud.ch"
$ws.Hyperlinks.Add($ws.Range("C21"), "mailto:tbd@eduvaud.ch")
$ws.Range("C3").Copy($ws.Range("C21"))

$ws.Range("D21").Value = "tbd@eduvaud.ch"
$ws.Hyperlinks.Add($ws.Range("D21"), "mailto:tbd@eduvaud.ch")
$ws.Range("D3").Copy($ws.Range("D21"))

$ws.Range("E21").Value = "eleve"
$ws.Range("F21").Value = "cin1c"

$ws.Range("G21").Value = 44774
$ws.Range("G3").Copy($ws.Range("G21"))

$ws.Range("H21").Value = "ruppture"

# --- Row 22: continuation log line for the "tbd" user ---
$ws.Range("C3").Copy($ws.Range("C22"))

$ws.Range("D22").Value = "tbd@eduvaud.ch"
$ws.Hyperlinks.Add($ws.Range("D22"), "mailto:tbd@eduvaud.ch")
$ws.Range("D3").Copy($ws.Range("D22"))

$ws.Range("G22").Value = 44774
$ws.Range("G3").Copy($ws.Range("G22"))

# --- Row 24: new "ghost" (orphaned) user, full sync-log line ---
$ws.Range("A24").Value = "ghost"
$ws.Range("B24").Value = "ghost"

$ws.Range("C24").Value = "ghost@eduvaud.ch"
$ws.Hyperlinks.Add($ws.Range("C24"), "mailto:ghost@eduvaud.ch")
$ws.Range("C3").Copy($ws.Range("C24"))

$ws.Range("D24").Value = "ghost@eduvaud.ch"
$ws.Hyperlinks.Add($ws.Range("D24"), "mailto:ghost@eduvaud.ch")
$ws.Range("D3").Copy($ws.Range("D24"))

$ws.Range("G24").Value = 44774
$ws.Range("G3").Copy($ws.Range("G24"))

$ws.Range("H24").Value = "rupture"

# --- back to row 22: append the trailing comment now that "rupture" exists ---
$ws.Range("H22").Value = "bla ruPture au 748"

# --- Row 23: blank separator log line ---
$ws.Range("C3").Copy($ws.Range("C23"))
$ws.Range("D3").Copy($ws.Range("D23"))
$ws.Range("G3").Copy($ws.Range("G23"))

# --- Row 27: new "tbr" (to-be-restored) user, full sync-log line ---
$ws.Range("A27").Value = "tbr"
$ws.Range("B27").Value = "tbr"

$ws.Range("C27").Value = "tbr@eduvaud.ch"
$ws.Hyperlinks.Add($ws.Range("C27"), "mailto:tbr@eduvaud.ch")
$ws.Range("C3").Copy($ws.Range("C27"))

$ws.Range("D27").Value = "tbr@eduvaud.ch"
$ws.Hyperlinks.Add($ws.Range("D27"), "mailto:tbr@eduvaud.ch")
$ws.Range("D3").Copy($ws.Range("D27"))

$ws.Range("E27").Value = "prof"
$ws.Range("F27").Value = "min1"

$ws.Range("G27").Value = 44774
$ws.Range("G3").Copy($ws.Range("G27"))

# --- Row 28: continuation log line for the "tbr" user ---
$ws.Range("C3").Copy($ws.Range("C28"))

$ws.Range("D28").Value = "tbr@eduvaud.ch"
$ws.Hyperlinks.Add($ws.Range("D28"), "mailto:tbr@eduvaud.ch")
$ws.Range("D3").Copy($ws.Range("D28"))

$ws.Range("G28").Value = 44774
$ws.Range("G3").Copy($ws.Range("G28"))

$ws.Range("H28").Value = "rupture"

# --- Row 29: "tbr" user restored, duplicate full sync-log line ---
$ws.Range("A29").Value = "tbr"
$ws.Range("B29").Value = "tbr"

$ws.Range("C29").Value = "tbr@eduvaud.ch"
$ws.Hyperlinks.Add($ws.Range("C29"), "mailto:tbr@eduvaud.ch")
$ws.Range("C3").Copy($ws.Range("C29"))

$ws.Range("D29").Value = "tbr@eduvaud.ch"
$ws.Hyperlinks.Add($ws.Range("D29"), "mailto:tbr@eduvaud.ch")
$ws.Range("D3").Copy($ws.Range("D29"))

$ws.Range("E29").Value = "prof"
$ws.Range("F29").Value = "min1"

$ws.Range("G29").Value = 44774
$ws.Range("G3").Copy($ws.Range("G29"))

$ws.Range("H29").Select()
